$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the third row (d9261fb8...)
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 20:57:58"

# zh-cn sheet: Correspond Handoff/Handback Datetime for d9261fb8... row
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-31 20:57:53"
$wsZhCn.Range("K4").Value = "2016-08-31 20:58:26"

# de-de sheet: Correspond Handoff/Handback Datetime for d9261fb8... row
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-31 20:57:58"
$wsDeDe.Range("K4").Value = "2016-08-31 20:58:34"
